$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.819.91'
$ws.Range('E2').Value = '  +2.92%  '

$ws.Range('D3').Value = '1.722.85'
$ws.Range('E3').Value = '  +2.85%  '

$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.54%  '

$ws.Range('D5').Value = '217.55'
$ws.Range('E5').Value = '  +1.05%  '

$ws.Range('D6').Value = '0.523'
$ws.Range('E6').Value = '  +1.16%  '

$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.61%  '

$ws.Range('D8').Value = '24.17'
$ws.Range('E8').Value = '  +12.77%  '

$ws.Range('E9').Value = '  +4.76%  '

$ws.Range('D10').Value = '0.0632'
$ws.Range('E10').Value = '  +1.48%  '

$ws.Range('D11').Value = '0.0898'
$ws.Range('E11').Value = '  +1.05%  '

$ws.Range('D12').Value = '1.967.37'
$ws.Range('E12').Value = '  +2.89%  '

$ws.Range('D13').Value = '1.727.24'
$ws.Range('E13').Value = '  +3.25%  '

$ws.Range('D14').Value = '4.25'
$ws.Range('E14').Value = '  +3.29%  '

$ws.Range('E15').Value = '  +5.95%  '

$ws.Range('D16').Value = '68.17'
$ws.Range('E16').Value = '  +2.77%  '

$ws.Range('D17').Value = '27.837.84'
$ws.Range('E17').Value = '  +2.96%  '

$ws.Range('D18').Value = '242.32'
$ws.Range('E18').Value = '  +2.63%  '

$ws.Range('D19').Value = '8.07'
$ws.Range('E19').Value = '  -1.19%  '

$ws.Range('D20').Value = '0.0₃0750'
$ws.Range('E20').Value = '  +1.91%  '

$ws.Range('D21').Value = '0.995'
$ws.Range('E21').Value = '  -0.60%  '

$ws.Range('D22').Value = '4.64'
$ws.Range('E22').Value = '  +3.94%  '

$ws.Range('E23').Value = '  +5.23%  '

$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('D25').Value = '148.73'
$ws.Range('E25').Value = '  +0.38%  '

$ws.Range('D26').Value = '7.56'
$ws.Range('E26').Value = '  +3.92%  '

$ws.Range('D27').Value = '16.65'
$ws.Range('E27').Value = '  +0.91%  '

$ws.Range('E28').Value = '  +1.07%  '

$ws.Range('D29').Value = '0.995'
$ws.Range('E29').Value = '  -0.61%  '

$ws.Range('D30').Value = '0.0507'
$ws.Range('E30').Value = '  +1.69%  '

$ws.Range('E31').Value = '  +1.19%  '

$ws.Range('E32').Value = '  +2.02%  '

$ws.Range('D33').Value = '1.544.95'
$ws.Range('E33').Value = '  +0.17%  '

$ws.Range('D34').Value = '3.32'
$ws.Range('E34').Value = '  +4.59%  '

$ws.Range('E35').Value = '  -1.44%  '

$ws.Range('E36').Value = '  +6.36%  '

$ws.Range('E37').Value = '  +4.41%  '

$ws.Range('E38').Value = '  +0.20%  '

$ws.Range('E39').Value = '  +0.05%  '

$ws.Range('E40').Value = '  +2.85%  '

$ws.Range('D41').Value = '71.51'
$ws.Range('E41').Value = '  +5.51%  '

$ws.Range('D42').Value = '5.87'
$ws.Range('E42').Value = '  +5.89%  '

$ws.Range('D43').Value = '0.995'
$ws.Range('E43').Value = '  -0.64%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.871.52'
$ws.Range('E44').Value = '  +2.81%  '

$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '2.28'
$ws.Range('E45').Value = '  +1.59%  '

$ws.Range('D46').Value = '0.790'
$ws.Range('E46').Value = '  +1.42%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '92.17'
$ws.Range('E47').Value = '  +1.90%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.69'
$ws.Range('E48').Value = '  +9.58%  '

$ws.Range('D49').Value = '0.0₆0109'
$ws.Range('E49').Value = '  +1.06%  '

$ws.Range('D50').Value = '8.37'
$ws.Range('E50').Value = '  +4.53%  '

$ws.Range('E51').Value = '  +1.99%  '
